{"js": "const ooxml = \"<?xml version=\\\"1.0\\\" encoding=\\\"UTF-8\\\" standalone=\\\"yes\\\"?><pkg:package xmlns:pkg=\\\"http://schemas.microsoft.com/office/2006/xmlPackage\\\"><pkg:part pkg:name=\\\"/word/document.xml\\\" pkg:contentType=\\\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\\\"><pkg:xmlData><w:document xmlns:w=\\\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\\\"><w:body><w:p><w:r><w:t>(?) denotes that sentence may not be 100% correct.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\\\"gramStart\\\"/><w:r><w:t>All of</w:t></w:r><w:proofErr w:type=\\\"gramEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> the drawing will have to be updated to work with the new unity graphics</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>JoyUpdate</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">, </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>JoyCalibrate</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">, and </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>JoyInit</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> are useless for what we are doing.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Silence appears to kill the sound and reset it (?)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>This game supports WASD, arrows, and numpad keys (4862).</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\\\"preserve\\\">There is an odd variable Elroy and </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>superElroy</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">. Is this </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>paku</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>paku</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t>?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\\\"preserve\\\">The </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>playEatGhost</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> procedure appears to reset the sound and play the eating ghost sound only, likely not required to be transferred</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>The scores in the pascal program uses longs for the value.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>On death the player is reset facing and moving to the left.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\\\"preserve\\\">Every movement is linked to an exit counter for each of the ghosts </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Blane and Hinky have special logics upon initialization. This appears to be linked to the exit counter to leave the jail in the center.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\\\"preserve\\\">We are likely to make the </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>menuLoop</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> procedure as our main and </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>gameLoop</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> as what it </w:t></w:r><w:proofErr w:type=\\\"gramStart\\\"/><w:r><w:t>calls, or</w:t></w:r><w:proofErr w:type=\\\"gramEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> base it off of that.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>gameLoop</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">: Displays the ready on each new level, then plays the starting theme on </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>firstRun</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> and sets </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>firstRun</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> to false, then will set the speed to the current level\\u2019s rate and </w:t></w:r><w:r><w:t xml:space=\\\"preserve\\\">resets the jail timers, </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>dotcounters</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">, and </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>fleetotal</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\">. Pellets seem to be on a loop of 6 for the blinking. Stinky also calls his update every three ticks of the timer(?) and then the game will check for ghost collisions. lines 2013 to 2079 are </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space=\\\"preserve\\\">sound handlers and can be ignored. Afterwards the game will check to see if a char got read in to change directions and repeat until the dots are &lt;=0 or lives are &lt;= 0. Dots = 0 game moves to next level. Lives &lt;= 0 or stage &lt;=0 the game checks to see if the </w:t></w:r><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>scorelist</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> can be updated and ends the procedure.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\\\"spellStart\\\"/><w:r><w:t>menuLoop</w:t></w:r><w:proofErr w:type=\\\"spellEnd\\\"/><w:r><w:t xml:space=\\\"preserve\\\"> controls the menu flashing and will go into the game if the key is correct, and will redraw the menu once the game is over and resume the menu loop.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\";\ncontext.document.body.insertOoxml(ooxml, Word.InsertLocation.start);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r><w:t>(?) denotes that sentence may not be 100% correct.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\"gramStart\"/><w:r><w:t>All of</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> the drawing will have to be updated to work with the new unity graphics</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>JoyUpdate</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>JoyCalibrate</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>JoyInit</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> are useless for what we are doing.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Silence appears to kill the sound and reset it (?)</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>This game supports WASD, arrows, and numpad keys (4862).</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\"preserve\">There is an odd variable Elroy and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>superElroy</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">. Is this </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>paku</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>paku</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>?</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\"preserve\">The </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>playEatGhost</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> procedure appears to reset the sound and play the eating ghost sound only, likely not required to be transferred</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>The scores in the pascal program uses longs for the value.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>On death the player is reset facing and moving to the left.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\"preserve\">Every movement is linked to an exit counter for each of the ghosts </w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Blane and Hinky have special logics upon initialization. This appears to be linked to the exit counter to leave the jail in the center.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space=\"preserve\">We are likely to make the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>menuLoop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> procedure as our main and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>gameLoop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> as what it </w:t></w:r><w:proofErr w:type=\"gramStart\"/><w:r><w:t>calls, or</w:t></w:r><w:proofErr w:type=\"gramEnd\"/><w:r><w:t xml:space=\"preserve\"> base it off of that.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>gameLoop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">: Displays the ready on each new level, then plays the starting theme on </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>firstRun</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> and sets </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>firstRun</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> to false, then will set the speed to the current level\u2019s rate and </w:t></w:r><w:r><w:t xml:space=\"preserve\">resets the jail timers, </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>dotcounters</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">, and </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>fleetotal</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\">. Pellets seem to be on a loop of 6 for the blinking. Stinky also calls his update every three ticks of the timer(?) and then the game will check for ghost collisions. lines 2013 to 2079 are </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space=\"preserve\">sound handlers and can be ignored. Afterwards the game will check to see if a char got read in to change directions and repeat until the dots are &lt;=0 or lives are &lt;= 0. Dots = 0 game moves to next level. Lives &lt;= 0 or stage &lt;=0 the game checks to see if the </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>scorelist</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> can be updated and ends the procedure.</w:t></w:r></w:p><w:p/><w:p><w:proofErr w:type=\"spellStart\"/><w:r><w:t>menuLoop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> controls the menu flashing and will go into the game if the key is correct, and will redraw the menu once the game is over and resume the menu loop.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n$r = $d.Range(0, 0)\n$r.InsertXML($xml)\n\n$paraCount = $d.Paragraphs.Count\n$secondLast = $d.Paragraphs($paraCount - 1)\n$mark = $d.Range($secondLast.Range.End - 1, $secondLast.Range.End)\n$mark.Delete()\n"}
